$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new column widths for the newly introduced columns C and D
$ws.Columns.Item(3).ColumnWidth = 57.5
$ws.Columns.Item(4).ColumnWidth = 44

# Row 1 (headers)
$ws.Range("C1").Value2 = "Default Coordinates"
$ws.Range("D1").Value2 = "Corresponding Brain Region"

# Row 2 - Film and TV
$ws.Range("C2").Value2 = "THREE.Vector3(-30.06, 4.40, -59.25)"
$ws.Range("D2").Value2 = "Occipital Lobe"

# Row 3 - Thoughts
$ws.Range("C3").Value2 = "THREE.Vector3(-15.55, 20.11, 67.40)"
$ws.Range("D3").Value2 = "Anterior prefrontal cortex"

# Row 4 - Books
$ws.Range("C4").Value2 = "THREE.Vector3(53.49, 15.30, -18.65)"
$ws.Range("D4").Value2 = "Temporal language areas (Wernicke’s area)"

# Row 5 - Music
$ws.Range("C5").Value2 = "THREE.Vector3(-55.86, 9.78, -3.42)"
$ws.Range("D5").Value2 = "Auditory cortex"

# Row 6 - Games and Interactive Media
$ws.Range("C6").Value2 = "THREE.Vector3(31.93, 45.64, 12.26)"
$ws.Range("D6").Value2 = "Motor cortex"

# Row 7 - Protocols
$ws.Range("C7").Value2 = "THREE.Vector3(-11.15, 47.78, 38.80)"
$ws.Range("D7").Value2 = "Dorsolateral prefrontal cortex"

# Row 8 - Principles
$ws.Range("C8").Value2 = "THREE.Vector3(38.62, 21.26, 40.53)"
$ws.Range("D8").Value2 = "Ventromedial prefrontal cortex"

# Row 9 - Food
$ws.Range("C9").Value2 = "THREE.Vector3(-15.16, -28.65, 1.76)"
$ws.Range("D9").Value2 = "Hypothalamus"

# Row 10 - Adventures
$ws.Range("C10").Value2 = "THREE.Vector3(-48.33, -31.49, -23.14)"
$ws.Range("D10").Value2 = "Cerebellum"

# Update the active selection to match the authored state
$ws.Range("C14").Select()
